# Scheduled market-data refresh: update cached currentAveragePrice /
# LevePrice / LeveProfit columns (H, I, J, K, L, M, N) per-leve across the
# job sheets with freshly pulled marketboard data. Cell coordinates below
# were located once per sheet via Items's Leve Item ID (column G) anchors.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 14869.5
$ws.Range("I76").Value = 14869.5
$ws.Range("K76").Value = 14869.5
$ws.Range("M76").Value = -14554.5
$ws.Range("H79").Value = 14869.5
$ws.Range("I79").Value = 14869.5
$ws.Range("K79").Value = 14869.5
$ws.Range("M79").Value = -13777.5
$ws.Range("H94").Value = 2699
$ws.Range("I94").Value = 2699
$ws.Range("K94").Value = 2699
$ws.Range("M94").Value = -2248
$ws.Range("H100").Value = 6830.231
$ws.Range("I100").Value = 5399.5713
$ws.Range("K100").Value = 5399.5713
$ws.Range("M100").Value = -4858.5713
$ws.Range("H131").Value = 3790369.2
$ws.Range("I131").Value = 2988.6
$ws.Range("K131").Value = 8965.799999999999
$ws.Range("M131").Value = -3925.799999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6524.296
$ws.Range("I32").Value = 5427.183
$ws.Range("K32").Value = 5427.183
$ws.Range("M32").Value = -5140.183
$ws.Range("H45").Value = 7942342.5
$ws.Range("J45").Value = 7748.5
$ws.Range("L45").Value = 7748.5
$ws.Range("N45").Value = -8502.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 741441.3
$ws.Range("I105").Value = 1145638.6
$ws.Range("K105").Value = 1145638.6
$ws.Range("M105").Value = -1143891.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 22731048
$ws.Range("I31").Value = 38464250
$ws.Range("K31").Value = 38464250
$ws.Range("M31").Value = -38463955
$ws.Range("H34").Value = 22731048
$ws.Range("I34").Value = 38464250
$ws.Range("K34").Value = 38464250
$ws.Range("M34").Value = -38464048
$ws.Range("H86").Value = 6530.231
$ws.Range("I86").Value = 6093.5
$ws.Range("K86").Value = 6093.5
$ws.Range("M86").Value = -4970.5
$ws.Range("H89").Value = 6530.231
$ws.Range("I89").Value = 6093.5
$ws.Range("K89").Value = 30467.5
$ws.Range("M89").Value = -24851.5
$ws.Range("H105").Value = 2708.1428
$ws.Range("I105").Value = 2708.1428
$ws.Range("K105").Value = 2708.1428
$ws.Range("M105").Value = -961.1428000000001
$ws.Range("H122").Value = 3618.75
$ws.Range("J122").Value = 1487.5
$ws.Range("L122").Value = 4462.5
$ws.Range("N122").Value = -9362.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 831.8095
$ws.Range("I97").Value = 740.4211
$ws.Range("K97").Value = 740.4211
$ws.Range("M97").Value = -244.4211
$ws.Range("H122").Value = 4715372.5
$ws.Range("I122").Value = 8250951.5
$ws.Range("K122").Value = 24752854.5
$ws.Range("M122").Value = -24750404.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4634.346
$ws.Range("I40").Value = 4630.5654
$ws.Range("K40").Value = 4630.5654
$ws.Range("M40").Value = -4494.5654
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = ""
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = ""
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").Value = ""
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").Value = ""
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").Value = ""
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").Value = ""
$ws.Range("H93").Value = 2927019.8
$ws.Range("J93").Value = 6949918
$ws.Range("L93").Value = 6949918
$ws.Range("N93").Value = -6952414
$ws.Range("H97").Value = 64344
$ws.Range("J97").Value = 64344
$ws.Range("L97").Value = 64344
$ws.Range("N97").Value = -66326
$ws.Range("H98").Value = 99996.5
$ws.Range("J98").Value = 99996.5
$ws.Range("L98").Value = 99996.5
$ws.Range("N98").Value = -105986.5
$ws.Range("H101").Value = 31648.834
$ws.Range("J101").Value = 31648.834
$ws.Range("L101").Value = 31648.834
$ws.Range("N101").Value = -38138.834
$ws.Range("H103").Value = 68551
$ws.Range("J103").Value = 68551
$ws.Range("L103").Value = 68551
$ws.Range("N103").Value = -70895
$ws.Range("H104").Value = 67282.5
$ws.Range("J104").Value = 67282.5
$ws.Range("L104").Value = 67282.5
$ws.Range("N104").Value = -74270.5
$ws.Range("H106").Value = 18000
$ws.Range("J106").Value = 18000
$ws.Range("L106").Value = 18000
$ws.Range("N106").Value = -20524
$ws.Range("H110").Value = 99998
$ws.Range("J110").Value = 99998
$ws.Range("L110").Value = 99998
$ws.Range("N110").Value = -108178
$ws.Range("H114").Value = 115449
$ws.Range("J114").Value = 115449
$ws.Range("L114").Value = 115449
$ws.Range("N114").Value = -124127
$ws.Range("H129").Value = 92021.5
$ws.Range("J129").Value = 92021.5
$ws.Range("L129").Value = 92021.5
$ws.Range("N129").Value = -102021.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 39891.332
$ws.Range("J70").Value = 39887
$ws.Range("L70").Value = 39887
$ws.Range("N70").Value = -40517
$ws.Range("H73").Value = 39891.332
$ws.Range("J73").Value = 39887
$ws.Range("L73").Value = 39887
$ws.Range("N73").Value = -42071
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").Value = ""
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").Value = ""
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").Value = ""
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").Value = ""
$ws.Range("H104").Value = 112650
$ws.Range("J104").Value = 112650
$ws.Range("L104").Value = 112650
$ws.Range("N104").Value = -119638
$ws.Range("H113").Value = 484.8387
$ws.Range("I113").Value = 404.45456
$ws.Range("K113").Value = 1213.36368
$ws.Range("M113").Value = 956.6363200000001
$ws.Range("H136").Value = 1282.1875
$ws.Range("I136").Value = 651.4737
$ws.Range("J136").Value = 2204
$ws.Range("K136").Value = 1954.4211
$ws.Range("L136").Value = 6612
$ws.Range("M136").Value = 595.5789
$ws.Range("N136").Value = -11712
